{"js": "// Bump the Zombie bite's IntrinsicWeapon damage value from 10 to 15 in the\n// design-rationale text: \"...return a new IntrinsicWeapon(10, \"bites\") instead\n// of the default...\" -> \"...IntrinsicWeapon(15, \"bites\") instead of the default...\".\nconst body = context.document.body;\n\nconst results = body.search(\"10, \\u201cbites\\u201d) instead of the default\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  target.insertText(\"15, \\u201cbites\\u201d) instead of the default\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Bump the Zombie bite's IntrinsicWeapon damage value from 10 to 15 in the\n# design-rationale text: \"...return a new IntrinsicWeapon(10, \"bites\") instead\n# of the default...\" -> \"...IntrinsicWeapon(15, \"bites\") instead of the default...\".\n$d = $word.ActiveDocument\n\n$openQuote = [char]0x201C\n$closeQuote = [char]0x201D\n\n$searchText = \"10, \" + $openQuote + \"bites\" + $closeQuote + \") instead of the default\"\n$replaceText = \"15, \" + $openQuote + \"bites\" + $closeQuote + \") instead of the default\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $searchText\n$find.Replacement.Text = $replaceText\n$find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n"}
